$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.814.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.292.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.889"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.636.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.289.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.721.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0918"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0383"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  +5.22%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +37.48%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
